$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 7")
$ws.Activate()

# Row 4: shorten the existing "V27, 28" note to just "V27" (the "28" part
# moves out into its own entries on rows 5 and 6 below).
$ws.Range("H4").Value = "V27"

# Row 5: new log entry (20:00-20:30, 30 min, activity "video", comment "V28", C=12)
$ws.Range("C5").Value = 0.83333333333333337
$ws.Range("D5").Value = 0.85416666666666663
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = "video"
$ws.Range("H5").Value = "V28"
$ws.Range("J5").Value = 12

# Row 6: new log entry (date 15/03/2020, start 14:15, activity "video", comment "V28")
$ws.Range("B6").Value = 43905
$ws.Range("C6").Value = 0.59375
$ws.Range("G6").Value = "video"
$ws.Range("H6").Value = "V28"

# Rows 7-10: activity column filled in with "video"
$ws.Range("G7").Value = "video"
$ws.Range("G8").Value = "video"
$ws.Range("G9").Value = "video"
$ws.Range("G10").Value = "video"

# Restore selection to the last-edited cell
$ws.Range("H6").Select()

$wb.Save()
